$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case: raw "string of 0/1s" in column A, with the computed
# "consecutive ones run-lengths (>1)" answer in column B.
$ws.Range("A11").Value = "1,1,0,0,0,0,0,0,1,1,1,1"

# The MAP/LAMBDA array formula anchored at B3 previously covered A3:A10 /
# B3:B10. Re-enter it (Ctrl+Shift+Enter style) over the full B3:B11 range so
# it now covers A3:A11 and spills its answer into the new row too.
$ws.Range("B3:B11").FormulaArray = "=MAP(A3:A11,LAMBDA(x,LET(a,LEN(SUBSTITUTE(TEXTSPLIT(x,0),`",`",`"`")),TEXTJOIN(`",`",1,IF(a>1,a,`"`")))))"

# Move the active selection to match where the author ended up.
$ws.Range("H5").Select()
